$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty cells D18:F18 with the value 5, matching
# the pattern already present in neighboring rows. This also updates the
# shared SUM formula result in J18 from 5 to 20.
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 5

# Update the active selection in the frozen bottom-right pane to G18,
# matching the author's last cursor position after editing the row.
$ws.Range("G18").Select() | Out-Null
